# Generate Report for Handoff
# Mark the "859bec39-cd0a-42b5-afc0-3f20e77eff8e" file row as "Ready for handoff"
# across the Overview, zh-cn and de-de worksheets, and stamp the new handoff
# timestamps for each locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 18:17:45"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 18:17:41"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 18:17:45"
